# Nouveau diagramme de gantt:
# - Extend the Gantt schedule with a new "Maintenace" step (row 9)
# - Recompute/extend downstream durations+dates on existing steps
# - Turn the A1:D9 range into an Excel Table ("Tableau1")
# - Apply a dd/mm/yyyy date format on the "Date de debut" / "Date de fin" columns

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-8 (durations in column C; formulas in B/D recalc automatically) ---
$ws.Range("C2").Value = 5
$ws.Range("C3").Value = 3
$ws.Range("C4").Value = 18
$ws.Range("C5").Value = 14
$ws.Range("C6").Value = 5
$ws.Range("C7").Value = 1
$ws.Range("C8").Value = 1

# --- Add the new "Maintenace" row 9 ---
$ws.Range("A9").Value = "Maintenace"
$ws.Range("B9").Formula = "=+D8+1"
$ws.Range("C9").Value = 30
$ws.Range("D9").Formula = "=+C9+B9"

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 15.7109375
$ws.Columns.Item(3).ColumnWidth = 16.7109375
$ws.Columns.Item(4).ColumnWidth = 12.85546875

# --- Selected cell moved to C7 in the source workbook ---
$ws.Range("C7").Select()

# --- Convert the data range into a proper Excel Table ---
$ws.Range("B2:B9").NumberFormat = "dd/mm/yyyy"
$ws.Range("D2:D9").NumberFormat = "dd/mm/yyyy"
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:D9"), $null, 1)
$lo.Name = "Tableau1"

# --- Update the Gantt chart series/category ranges & axis scaling ---
$chart = $ws.ChartObjects(1).Chart
$chart.SeriesCollection(1).XValues = $ws.Range("A2:A9")
$chart.SeriesCollection(1).Values = $ws.Range("B2:B9")
$chart.SeriesCollection(2).Values = $ws.Range("C2:C9")
$chart.Axes(2).MaximumScale = 45165
